$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.121.39'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.175.16'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.611'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.08'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -6.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.69'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '36.41'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -10.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0934'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.84'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.500.81'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.24'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.843'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.187.97'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.036.78'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0943'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.40'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.03'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.46'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.82'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.24%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.39'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.61'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.04'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.10'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.79%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.65'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0744'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.51'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.68'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0303'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +16.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.20'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.96%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '60.44'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.23'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -11.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.48'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.188'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -7.20%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0988'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.27'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -10.12%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.08%  '
